# Apply updated crypto price/volume data per commit "Updated cryptos list on Sat Sep 23 23:41:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'26.717.51"
$ws.Cells.Item(2, 5).Value = "  +0.22%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.600.59"
$ws.Cells.Item(3, 5).Value = "  +0.23%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.31%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'211.61"
$ws.Cells.Item(5, 5).Value = "  -0.07%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.67%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.30%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.11%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.84%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'19.57"
$ws.Cells.Item(10, 5).Value = "  +0.18%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.0842"
$ws.Cells.Item(11, 5).Value = "  +0.65%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'1.824.74"
$ws.Cells.Item(12, 5).Value = "  +0.21%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'1.607.86"
$ws.Cells.Item(13, 5).Value = "  -0.41%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.53%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +0.17%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'65.36"
$ws.Cells.Item(16, 5).Value = "  +1.41%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'26.691.09"
$ws.Cells.Item(17, 5).Value = "  +0.22%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +2.98%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +3.77%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.27%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'209.10"
$ws.Cells.Item(21, 5).Value = "  +0.05%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.51%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.70%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'8.95"
$ws.Cells.Item(24, 5).Value = "  +0.62%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'142.49"
$ws.Cells.Item(25, 5).Value = "  -1.93%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.35%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.83%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.15%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.54%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'0.0522"
$ws.Cells.Item(30, 5).Value = "  +2.96%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'1.16"
$ws.Cells.Item(31, 5).Value = "  -0.41%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.26"

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.62%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.294.25"
$ws.Cells.Item(34, 5).Value = "  +1.30%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.625"
$ws.Cells.Item(35, 5).Value = "  -5.20%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.90%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.24%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +20.21%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -2.33%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.82%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.04%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.49%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'63.19"
$ws.Cells.Item(44, 5).Value = "  -2.16%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'1.736.66"
$ws.Cells.Item(45, 5).Value = "  +0.13%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'91.35"
$ws.Cells.Item(46, 5).Value = "  +1.46%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -1.85%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).Value = "'0.101"
$ws.Cells.Item(48, 5).Value = "  -1.37%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).Value = "'0.0510"
$ws.Cells.Item(49, 5).Value = "  +0.65%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "USDD"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(50, 4).Value = "'1.00"
$ws.Cells.Item(50, 5).Value = "  +0.30%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "'7.41"
$ws.Cells.Item(51, 5).Value = "  -0.60%  "
